$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the "Apio" price series. In the source
# data the records are ordered chronologically-ish but this new one lands
# at the top of the tail block (row 257), pushing the existing rows
# 257:271 down to 258:272 (dimension grows from R271 to R272).
$ws.Rows.Item(257).Insert()

# Populate the freshly inserted row 257 with the new record. All of the
# "fixed" columns (market/region/category/etc.) are identical to the rest
# of this block; only the date (Fecha) and the associated price/volume
# figures differ.
$ws.Range("A257").Value = 5
$ws.Range("B257").Value = "Macroferia Regional de Talca"
$ws.Range("C257").Value = "Maule"
$ws.Range("D257").Value = 45021
$ws.Range("E257").Value = 7
$ws.Range("F257").Value = 100112017
$ws.Range("G257").Value = "Apio"
$ws.Range("H257").Value = "Americana (o)"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 500
$ws.Range("K257").Value = 8000
$ws.Range("L257").Value = 8000
$ws.Range("M257").Value = 8000
$ws.Range("N257").Value = "`$/docena de matas"
$ws.Range("O257").Value = "Provincia del Elquí"
$ws.Range("P257").Value = 1333
$ws.Range("Q257").Value = 6
$ws.Range("R257").Value = "Hortaliza"
